$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's mod-count entry as a new row (row 32)
# Format A32 as text first so the date-like string "2025/12/11" is stored
# as literal text, matching the existing Date column entries, instead of
# being auto-converted into a date serial number.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2025/12/11"
$ws.Range("B32").Value = "逃离鸭科夫"
$ws.Range("C32").Value = 1354

# Match the formatting of the preceding data rows (centered alignment,
# General number format) by copying row 31's formatting onto row 32.
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)  # xlPasteFormats
